$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Clear existing Item Name / UOM cells so their shared strings become unused
# and get dropped, letting us rebuild the shared-string table in the desired order.
for ($r = 2; $r -le 22; $r++) {
    $ws.Cells.Item($r, 4).Value = ""
    $ws.Cells.Item($r, 5).Value = ""
}

# Set Item Name (column D) for all rows first, in row order,
# then UOM (column E) for all rows, to reproduce shared-string build order.
$ws.Cells.Item(2, 4).Value = 'Esoral 20mg Capsule 50''s'
$ws.Cells.Item(3, 4).Value = 'Esoral 20mg Tablet  80''s'
$ws.Cells.Item(4, 4).Value = 'Esoral Injection & Capsule 20'
$ws.Cells.Item(5, 4).Value = 'Esoral 20mg Capsule Container 30''s'
$ws.Cells.Item(6, 4).Value = 'Esoral 20mg Tablet - 50''s'
$ws.Cells.Item(7, 4).Value = 'Esoral 40mg Tablet'
$ws.Cells.Item(8, 4).Value = 'Esoral 40mg EC Tablet - 42''s'
$ws.Cells.Item(9, 4).Value = 'Esoral 20mg Tablet'
$ws.Cells.Item(10, 4).Value = 'Esoral 20mg  Tablet 100''s'
$ws.Cells.Item(11, 4).Value = 'Esoral Injection & MUPS 20'
$ws.Cells.Item(12, 4).Value = 'Losectil 40mg Capsule (24''s)'
$ws.Cells.Item(13, 4).Value = 'Losectil 20mg Powder for Oral Suspension - 30''s'
$ws.Cells.Item(14, 4).Value = 'Losectil 10mg Capsule'
$ws.Cells.Item(15, 4).Value = 'Losectil DR Tablet '
$ws.Cells.Item(16, 4).Value = 'Losectil 40mg Powder for Oral Suspension'
$ws.Cells.Item(17, 4).Value = 'Losectil 20mg Powder for Oral Suspension'
$ws.Cells.Item(18, 4).Value = 'Losectil 40mg Capsule - 48''s'
$ws.Cells.Item(19, 4).Value = 'Losectil 20mg Capsule (100''s)'
$ws.Cells.Item(20, 4).Value = 'Losectil 20mg Capsule 500s'
$ws.Cells.Item(21, 4).Value = 'Rabifast 20mg Tablet - 50''s'
$ws.Cells.Item(22, 4).Value = 'Softi Ointment 15gm'

$ws.Cells.Item(2, 5).Value = '50''s'
$ws.Cells.Item(3, 5).Value = '80''s'
$ws.Cells.Item(4, 5).Value = 'Bundle'
$ws.Cells.Item(5, 5).Value = '30''s'
$ws.Cells.Item(6, 5).Value = '50''s'
$ws.Cells.Item(7, 5).Value = '20''s'
$ws.Cells.Item(8, 5).Value = '42''s'
$ws.Cells.Item(9, 5).Value = '20''s'
$ws.Cells.Item(10, 5).Value = '100''s'
$ws.Cells.Item(11, 5).Value = 'Bundle'
$ws.Cells.Item(12, 5).Value = '24 ''s'
$ws.Cells.Item(13, 5).Value = '30''s'
$ws.Cells.Item(14, 5).Value = '48 ''s'
$ws.Cells.Item(15, 5).Value = '60 ''s'
$ws.Cells.Item(16, 5).Value = '20''s'
$ws.Cells.Item(17, 5).Value = '20''s'
$ws.Cells.Item(18, 5).Value = '48''s'
$ws.Cells.Item(19, 5).Value = '100 ''s'
$ws.Cells.Item(20, 5).Value = '500''s'
$ws.Cells.Item(21, 5).Value = '50''s'
$ws.Cells.Item(22, 5).Value = '15gm'

# Update BSL NO (column A) per row
$ws.Cells.Item(2, 1).Value = 53
$ws.Cells.Item(3, 1).Value = 53
$ws.Cells.Item(4, 1).Value = 53
$ws.Cells.Item(5, 1).Value = 53
$ws.Cells.Item(6, 1).Value = 53
$ws.Cells.Item(7, 1).Value = 53
$ws.Cells.Item(8, 1).Value = 53
$ws.Cells.Item(9, 1).Value = 53
$ws.Cells.Item(10, 1).Value = 53
$ws.Cells.Item(11, 1).Value = 53
$ws.Cells.Item(12, 1).Value = 96
$ws.Cells.Item(13, 1).Value = 96
$ws.Cells.Item(14, 1).Value = 96
$ws.Cells.Item(15, 1).Value = 96
$ws.Cells.Item(16, 1).Value = 96
$ws.Cells.Item(17, 1).Value = 96
$ws.Cells.Item(18, 1).Value = 96
$ws.Cells.Item(19, 1).Value = 96
$ws.Cells.Item(20, 1).Value = 96
$ws.Cells.Item(21, 1).Value = 143
$ws.Cells.Item(22, 1).Value = 165
